$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the parenthetical abbreviation suffix from the "Name" column (C)
# descriptions for the monosaccharide residue rows. Row 10 (Sia / "Sialic
# acid") is intentionally left untouched - its text never had a
# parenthetical suffix.
$ws.Range("C2").Value  = "Hexose"
$ws.Range("C3").Value  = "Mannose"
$ws.Range("C4").Value  = "Galactose"
$ws.Range("C5").Value  = "Glucose"
$ws.Range("C6").Value  = "N-acetyl Hexosamine"
$ws.Range("C7").Value  = "N-acetyl Mannosamine"
$ws.Range("C8").Value  = "N-acetly Glucosamine"
$ws.Range("C9").Value  = "N-acetyl Galactosamine"
$ws.Range("C11").Value = "N-acetyl Neuraminic acid"
$ws.Range("C12").Value = "N-glycolyl Neuraminic acid"
$ws.Range("C13").Value = "Keto-Deoxy-Nonulonic acid"
$ws.Range("C14").Value = "deoxy Hexose"
$ws.Range("C15").Value = "Fucose"
$ws.Range("C16").Value = "Hexuronic acid"
$ws.Range("C17").Value = "Glucuronic acid"
$ws.Range("C18").Value = "Galacturonic acid"
$ws.Range("C19").Value = "Iduronic acid"
$ws.Range("C20").Value = "Mannuronic acid"
$ws.Range("C21").Value = "Hexosamine"
$ws.Range("C22").Value = "Glucosamine"
$ws.Range("C23").Value = "Galactosamine"
$ws.Range("C24").Value = "Mannosamine"
$ws.Range("C25").Value = "Pentose"
$ws.Range("C26").Value = "Xylose"

# Move the active selection from J23 to K15
$ws.Range("K15").Select() | Out-Null
